# Golden Valley Pipeline project documentation workbook update
# - Add a "py_files_doc" header cell ("Main Folder") with its column width
# - Move the active/selected tab from "column_info" to "changelog"

$wb = $excel.ActiveWorkbook

# --- py_files_doc sheet: add header value, size the column, and set selection ---
$wsPyFiles = $wb.Worksheets.Item("py_files_doc")
$wsPyFiles.Range("A1").Value = "Main Folder"
$wsPyFiles.Columns.Item(1).ColumnWidth = 17
$wsPyFiles.Range("A2").Select()

# --- changelog sheet becomes the active tab (was column_info) ---
$wsChangelog = $wb.Worksheets.Item("changelog")
$wsChangelog.Select()
